# Refresh the cryptocurrency "Price" and "Volume(1h)" columns (D, E) with
# the latest scraped figures, mirroring the scheduled GitHub Actions update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some of the new Price values (e.g. "1.00", "79.00") look like plain numbers,
# so force those specific cells to Text format first -- otherwise Excel would
# silently convert them (e.g. "1.00" -> 1) and the trailing zeros would be
# lost, unlike the source data which stores prices as text.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row -> (new Price, new Volume(1h)) updates
$ws.Cells.Item(2, 4).Value = "69.620.88"
$ws.Cells.Item(2, 5).Value = "  +5.02%  "
$ws.Cells.Item(3, 4).Value = "3.615.75"
$ws.Cells.Item(3, 5).Value = "  +5.16%  "
$ws.Cells.Item(4, 4).Value = "1.00"
$ws.Cells.Item(4, 5).Value = "  -0.11%  "
$ws.Cells.Item(5, 4).Value = "631.84"
$ws.Cells.Item(5, 5).Value = "  +5.69%  "
$ws.Cells.Item(6, 4).Value = "158.86"
$ws.Cells.Item(6, 5).Value = "  +8.66%  "
$ws.Cells.Item(7, 4).Value = "3.613.40"
$ws.Cells.Item(7, 5).Value = "  +5.10%  "
$ws.Cells.Item(8, 5).Value = "  -0.15%  "
$ws.Cells.Item(9, 5).Value = "  +4.70%  "
$ws.Cells.Item(10, 5).Value = "  +11.53%  "
$ws.Cells.Item(11, 4).Value = "7.53"
$ws.Cells.Item(11, 5).Value = "  +10.29%  "
$ws.Cells.Item(12, 5).Value = "  +7.12%  "
$ws.Cells.Item(13, 5).Value = "  +7.34%  "
$ws.Cells.Item(14, 4).Value = "33.78"
$ws.Cells.Item(14, 5).Value = "  +10.14%  "
$ws.Cells.Item(15, 4).Value = "4.225.53"
$ws.Cells.Item(15, 5).Value = "  +4.95%  "
$ws.Cells.Item(16, 4).Value = "3.611.46"
$ws.Cells.Item(16, 5).Value = "  +4.99%  "
$ws.Cells.Item(17, 4).Value = "69.543.99"
$ws.Cells.Item(17, 5).Value = "  +4.78%  "
$ws.Cells.Item(18, 5).Value = "  +1.33%  "
$ws.Cells.Item(19, 5).Value = "  +7.67%  "
$ws.Cells.Item(20, 4).Value = "16.21"
$ws.Cells.Item(20, 5).Value = "  +10.36%  "
$ws.Cells.Item(21, 4).Value = "10.31"
$ws.Cells.Item(21, 5).Value = "  +16.29%  "
$ws.Cells.Item(22, 4).Value = "463.14"
$ws.Cells.Item(22, 5).Value = "  +6.44%  "
$ws.Cells.Item(23, 4).Value = "0.648"
$ws.Cells.Item(23, 5).Value = "  +5.52%  "
$ws.Cells.Item(24, 4).Value = "79.00"
$ws.Cells.Item(24, 5).Value = "  +3.41%  "
$ws.Cells.Item(25, 5).Value = "  +12.17%  "
$ws.Cells.Item(26, 4).Value = "10.79"
$ws.Cells.Item(26, 5).Value = "  +9.18%  "
$ws.Cells.Item(27, 4).Value = "3.758.27"
$ws.Cells.Item(27, 5).Value = "  +4.88%  "
$ws.Cells.Item(28, 5).Value = "  +0.02%  "
$ws.Cells.Item(29, 4).Value = "9.45"
$ws.Cells.Item(29, 5).Value = "  +16.54%  "
$ws.Cells.Item(30, 4).Value = "2.66"
$ws.Cells.Item(30, 5).Value = "  +7.34%  "
$ws.Cells.Item(31, 4).Value = "1.74"
$ws.Cells.Item(31, 5).Value = "  +14.94%  "
$ws.Cells.Item(32, 5).Value = "  +8.84%  "
$ws.Cells.Item(33, 4).Value = "6.59"
$ws.Cells.Item(33, 5).Value = "  +9.46%  "
$ws.Cells.Item(34, 5).Value = "  +0.01%  "
$ws.Cells.Item(35, 5).Value = "  +7.85%  "
$ws.Cells.Item(36, 4).Value = "26.60"
$ws.Cells.Item(36, 5).Value = "  +5.61%  "
$ws.Cells.Item(37, 4).Value = "3.611.72"
$ws.Cells.Item(37, 5).Value = "  +5.31%  "
$ws.Cells.Item(38, 5).Value = "  +9.23%  "
$ws.Cells.Item(39, 4).Value = "2.42"
$ws.Cells.Item(39, 5).Value = "  +15.96%  "
$ws.Cells.Item(40, 5).Value = "  -0.01%  "
$ws.Cells.Item(41, 4).Value = "0.0930"
$ws.Cells.Item(41, 5).Value = "  +9.50%  "
$ws.Cells.Item(42, 4).Value = "178.91"
$ws.Cells.Item(42, 5).Value = "  +3.59%  "
$ws.Cells.Item(43, 4).Value = "1.00"
$ws.Cells.Item(43, 5).Value = "  +0.01%  "
$ws.Cells.Item(44, 4).Value = "5.71"
$ws.Cells.Item(44, 5).Value = "  +7.50%  "
$ws.Cells.Item(45, 4).Value = "32.13"
$ws.Cells.Item(45, 5).Value = "  +25.68%  "
$ws.Cells.Item(46, 4).Value = "0.915"
$ws.Cells.Item(46, 5).Value = "  +5.22%  "
$ws.Cells.Item(47, 4).Value = "1.40"
$ws.Cells.Item(47, 5).Value = "  +16.38%  "
$ws.Cells.Item(48, 4).Value = "2.77"
$ws.Cells.Item(48, 5).Value = "  +14.17%  "
$ws.Cells.Item(49, 4).Value = "45.98"
$ws.Cells.Item(49, 5).Value = "  +1.72%  "
$ws.Cells.Item(50, 5).Value = "  +5.56%  "
$ws.Cells.Item(51, 4).Value = "0.270"
$ws.Cells.Item(51, 5).Value = "  +12.18%  "
